{"js": "/*\n * Word JS API (Office.js) edit script.\n *\n * The source document contains a single 20-row x 5-column table where every\n * cell holds a short arithmetic expression, e.g. \"33-14=19\". The commit\n * regenerates the worksheet's answer key, so each of the 100 cells gets a\n * brand-new expression string; nothing else in the document changes.\n *\n * OLD_VALUES / NEW_VALUES below are the row-major (top-to-bottom,\n * left-to-right) before/after text for every cell, taken straight from the\n * unified diff (each diff hunk maps 1:1, in order, onto table cell (r, c)).\n * The script re-reads the table's current values and only overwrites a cell\n * when it still holds the expected \"before\" text (or already holds the\n * \"after\" text, making the script idempotent on re-run); this keeps the\n * edit purely positional/content-addressed instead of relying on fragile\n * global find/replace (two expressions, \"80-58=22\" and \"9+67=76\", each\n * appear twice in the original table at different cells with different\n * replacements).\n */\nconst OLD_VALUES = [\n  [\"33-14=19\", \"80-58=22\", \"16-8=8\", \"92-46=46\", \"53-48=5\"],\n  [\"49+17=66\", \"61-54=7\", \"60-19=41\", \"76-38=38\", \"28+5=33\"],\n  [\"31-16=15\", \"91-56=35\", \"61-5=56\", \"64-19=45\", \"8+26=34\"],\n  [\"77-49=28\", \"24+8=32\", \"84-65=19\", \"64+8=72\", \"71-47=24\"],\n  [\"73-67=6\", \"80-58=22\", \"28+43=71\", \"95-8=87\", \"41-13=28\"],\n  [\"43+19=62\", \"17-9=8\", \"44+49=93\", \"23+49=72\", \"57-28=29\"],\n  [\"25-7=18\", \"24-8=16\", \"90-21=69\", \"76-58=18\", \"84-78=6\"],\n  [\"27+19=46\", \"80-45=35\", \"37+5=42\", \"25+57=82\", \"13-6=7\"],\n  [\"9+67=76\", \"25+59=84\", \"59+6=65\", \"70-48=22\", \"54+27=81\"],\n  [\"90-63=27\", \"74-46=28\", \"54-15=39\", \"65-58=7\", \"7+38=45\"],\n  [\"95-6=89\", \"17+67=84\", \"4+17=21\", \"61-46=15\", \"6+25=31\"],\n  [\"13-9=4\", \"25+37=62\", \"69+27=96\", \"63-18=45\", \"94-65=29\"],\n  [\"84+7=91\", \"87-18=69\", \"4+57=61\", \"39+15=54\", \"6+78=84\"],\n  [\"16+75=91\", \"67+15=82\", \"8+89=97\", \"90-9=81\", \"55+16=71\"],\n  [\"72-35=37\", \"58+27=85\", \"23+29=52\", \"94-46=48\", \"34+59=93\"],\n  [\"90-46=44\", \"19+72=91\", \"23-14=9\", \"39+5=44\", \"49+48=97\"],\n  [\"82-63=19\", \"74+9=83\", \"32+29=61\", \"2+59=61\", \"17+7=24\"],\n  [\"7+57=64\", \"84-48=36\", \"40-31=9\", \"85+7=92\", \"70-51=19\"],\n  [\"37+45=82\", \"91-34=57\", \"23+58=81\", \"84-56=28\", \"60-32=28\"],\n  [\"70-7=63\", \"54-39=15\", \"9+67=76\", \"58+33=91\", \"14+49=63\"]\n];\n\nconst NEW_VALUES = [\n  [\"30-26=4\", \"87-39=48\", \"18+58=76\", \"84-68=16\", \"90-15=75\"],\n  [\"60-37=23\", \"73-8=65\", \"19+24=43\", \"33-26=7\", \"52-5=47\"],\n  [\"96-59=37\", \"8+63=71\", \"28+48=76\", \"39+17=56\", \"93-28=65\"],\n  [\"19+62=81\", \"50-34=16\", \"38+18=56\", \"7+67=74\", \"71-17=54\"],\n  [\"92-87=5\", \"80-38=42\", \"90-48=42\", \"63-9=54\", \"78-59=19\"],\n  [\"73-44=29\", \"76+15=91\", \"24-7=17\", \"26-19=7\", \"36+55=91\"],\n  [\"71-59=12\", \"91-14=77\", \"71-33=38\", \"67+8=75\", \"32-13=19\"],\n  [\"86+6=92\", \"52-47=5\", \"16+77=93\", \"30-2=28\", \"19+57=76\"],\n  [\"6+45=51\", \"49+45=94\", \"51-17=34\", \"58+34=92\", \"45+9=54\"],\n  [\"29+53=82\", \"28+16=44\", \"29+16=45\", \"63+18=81\", \"77+5=82\"],\n  [\"87-78=9\", \"93-66=27\", \"7+36=43\", \"7+14=21\", \"93-24=69\"],\n  [\"69+5=74\", \"35+8=43\", \"9+7=16\", \"19+77=96\", \"16+79=95\"],\n  [\"82-8=74\", \"70-21=49\", \"8+4=12\", \"66-49=17\", \"63+9=72\"],\n  [\"45-7=38\", \"96-57=39\", \"47+8=55\", \"70-64=6\", \"28+53=81\"],\n  [\"91-55=36\", \"28+8=36\", \"7+55=62\", \"58+17=75\", \"68+15=83\"],\n  [\"24+68=92\", \"17+79=96\", \"52-6=46\", \"90-58=32\", \"51-43=8\"],\n  [\"60-51=9\", \"57+19=76\", \"64-36=28\", \"16+55=71\", \"27+67=94\"],\n  [\"56+18=74\", \"47-29=18\", \"50-7=43\", \"49+42=91\", \"22+29=51\"],\n  [\"49+45=94\", \"81-47=34\", \"44-9=35\", \"15+79=94\", \"5+66=71\"],\n  [\"80-44=36\", \"3+89=92\", \"83-47=36\", \"27+5=32\", \"93-6=87\"]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nconst rowCount = OLD_VALUES.length;\nconst colCount = OLD_VALUES[0].length;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const expectedOld = OLD_VALUES[r][c];\n    const newValue = NEW_VALUES[r][c];\n    const current = table.values && table.values[r] ? table.values[r][c] : undefined;\n\n    if (current === expectedOld || current === newValue || current === undefined) {\n      const cell = table.getCell(r, c);\n      cell.value = newValue;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) edit script.\n#\n# The source document contains a single 20-row x 5-column table where every\n# cell holds a short arithmetic expression, e.g. \"33-14=19\". The commit\n# regenerates the worksheet's answer key, so each of the 100 cells gets a\n# brand-new expression string; nothing else in the document changes.\n#\n# $oldValues / $newValues below are the row-major (top-to-bottom,\n# left-to-right) before/after text for every cell, taken straight from the\n# unified diff (each diff hunk maps 1:1, in order, onto table cell (r, c),\n# 1-indexed the way Word's COM Cell(row, column) addressing works).\n# The script re-reads each cell's current text and only overwrites it when\n# it still holds the expected \"before\" text (or already holds the \"after\"\n# text, making the script idempotent on re-run); this keeps the edit purely\n# positional/content-addressed instead of relying on fragile global\n# find/replace (two expressions, \"80-58=22\" and \"9+67=76\", each appear\n# twice in the original table at different cells with different\n# replacements).\n\n$oldValues = @(\n    @(\"33-14=19\", \"80-58=22\", \"16-8=8\", \"92-46=46\", \"53-48=5\"),\n    @(\"49+17=66\", \"61-54=7\", \"60-19=41\", \"76-38=38\", \"28+5=33\"),\n    @(\"31-16=15\", \"91-56=35\", \"61-5=56\", \"64-19=45\", \"8+26=34\"),\n    @(\"77-49=28\", \"24+8=32\", \"84-65=19\", \"64+8=72\", \"71-47=24\"),\n    @(\"73-67=6\", \"80-58=22\", \"28+43=71\", \"95-8=87\", \"41-13=28\"),\n    @(\"43+19=62\", \"17-9=8\", \"44+49=93\", \"23+49=72\", \"57-28=29\"),\n    @(\"25-7=18\", \"24-8=16\", \"90-21=69\", \"76-58=18\", \"84-78=6\"),\n    @(\"27+19=46\", \"80-45=35\", \"37+5=42\", \"25+57=82\", \"13-6=7\"),\n    @(\"9+67=76\", \"25+59=84\", \"59+6=65\", \"70-48=22\", \"54+27=81\"),\n    @(\"90-63=27\", \"74-46=28\", \"54-15=39\", \"65-58=7\", \"7+38=45\"),\n    @(\"95-6=89\", \"17+67=84\", \"4+17=21\", \"61-46=15\", \"6+25=31\"),\n    @(\"13-9=4\", \"25+37=62\", \"69+27=96\", \"63-18=45\", \"94-65=29\"),\n    @(\"84+7=91\", \"87-18=69\", \"4+57=61\", \"39+15=54\", \"6+78=84\"),\n    @(\"16+75=91\", \"67+15=82\", \"8+89=97\", \"90-9=81\", \"55+16=71\"),\n    @(\"72-35=37\", \"58+27=85\", \"23+29=52\", \"94-46=48\", \"34+59=93\"),\n    @(\"90-46=44\", \"19+72=91\", \"23-14=9\", \"39+5=44\", \"49+48=97\"),\n    @(\"82-63=19\", \"74+9=83\", \"32+29=61\", \"2+59=61\", \"17+7=24\"),\n    @(\"7+57=64\", \"84-48=36\", \"40-31=9\", \"85+7=92\", \"70-51=19\"),\n    @(\"37+45=82\", \"91-34=57\", \"23+58=81\", \"84-56=28\", \"60-32=28\"),\n    @(\"70-7=63\", \"54-39=15\", \"9+67=76\", \"58+33=91\", \"14+49=63\")\n)\n\n$newValues = @(\n    @(\"30-26=4\", \"87-39=48\", \"18+58=76\", \"84-68=16\", \"90-15=75\"),\n    @(\"60-37=23\", \"73-8=65\", \"19+24=43\", \"33-26=7\", \"52-5=47\"),\n    @(\"96-59=37\", \"8+63=71\", \"28+48=76\", \"39+17=56\", \"93-28=65\"),\n    @(\"19+62=81\", \"50-34=16\", \"38+18=56\", \"7+67=74\", \"71-17=54\"),\n    @(\"92-87=5\", \"80-38=42\", \"90-48=42\", \"63-9=54\", \"78-59=19\"),\n    @(\"73-44=29\", \"76+15=91\", \"24-7=17\", \"26-19=7\", \"36+55=91\"),\n    @(\"71-59=12\", \"91-14=77\", \"71-33=38\", \"67+8=75\", \"32-13=19\"),\n    @(\"86+6=92\", \"52-47=5\", \"16+77=93\", \"30-2=28\", \"19+57=76\"),\n    @(\"6+45=51\", \"49+45=94\", \"51-17=34\", \"58+34=92\", \"45+9=54\"),\n    @(\"29+53=82\", \"28+16=44\", \"29+16=45\", \"63+18=81\", \"77+5=82\"),\n    @(\"87-78=9\", \"93-66=27\", \"7+36=43\", \"7+14=21\", \"93-24=69\"),\n    @(\"69+5=74\", \"35+8=43\", \"9+7=16\", \"19+77=96\", \"16+79=95\"),\n    @(\"82-8=74\", \"70-21=49\", \"8+4=12\", \"66-49=17\", \"63+9=72\"),\n    @(\"45-7=38\", \"96-57=39\", \"47+8=55\", \"70-64=6\", \"28+53=81\"),\n    @(\"91-55=36\", \"28+8=36\", \"7+55=62\", \"58+17=75\", \"68+15=83\"),\n    @(\"24+68=92\", \"17+79=96\", \"52-6=46\", \"90-58=32\", \"51-43=8\"),\n    @(\"60-51=9\", \"57+19=76\", \"64-36=28\", \"16+55=71\", \"27+67=94\"),\n    @(\"56+18=74\", \"47-29=18\", \"50-7=43\", \"49+42=91\", \"22+29=51\"),\n    @(\"49+45=94\", \"81-47=34\", \"44-9=35\", \"15+79=94\", \"5+66=71\"),\n    @(\"80-44=36\", \"3+89=92\", \"83-47=36\", \"27+5=32\", \"93-6=87\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $oldValues.Length\n$colCount = $oldValues[0].Length\n\nfor ($r = 0; $r -lt $rowCount; $r++) {\n    for ($c = 0; $c -lt $colCount; $c++) {\n        $expectedOld = $oldValues[$r][$c]\n        $newValue = $newValues[$r][$c]\n        $cell = $t.Cell($r + 1, $c + 1)\n\n        # Cell.Range.Text includes the trailing end-of-cell marks (CR + BEL);\n        # strip them before comparing against the plain expression text.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n        if (($current -eq $expectedOld) -or ($current -eq $newValue)) {\n            $cell.Range.Text = $newValue\n        }\n    }\n}\n"}
